$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, border, center/top alignment) from
# the existing header cell H1 onto the two new header cells I1 and J1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-6
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 10

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 6

$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 2

$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 8
